$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder / update existing rows: Water, Fire, Meat (with some value tweaks),
# then append Sun and Wheat.
$ws.Range("B2").Value = "Water"

$ws.Range("B3").Value = "Fire"
$ws.Range("F3").Value = 2

$ws.Range("E4").Value = 6
$ws.Range("H4").Value = 8

# New row 5: Sun
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Sun"
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 8
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 12
$ws.Range("I5").Value = 0.1
$ws.Range("J5").Value = 6

# New row 6: Wheat
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Wheat"
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 15
$ws.Range("I6").Value = 0.1
$ws.Range("J6").Value = 8

# Copy style from row 4 onto the two new rows so they match existing formatting
$ws.Range("A4:J4").Copy()
$ws.Range("A5:J6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("J7").Select()
